$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -1
    3  = 9
    4  = -2
    5  = -2
    6  = -3
    8  = 1
    9  = -1
    11 = 1
    12 = 9
    13 = -3
    14 = -4
    15 = 1
    16 = 2
    17 = -1
    18 = -3
    19 = 0
    20 = 4
    21 = -2
    22 = 5
    23 = -1
    24 = 8
    25 = -2
    26 = 6
    27 = 5
    28 = -1
    29 = -2
    30 = -1
    31 = -5
    32 = 1
    33 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
